$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - it is being removed from the dataset.
$ws.Rows(26).Delete()

# After the above deletion, the row that was "SC 92" (originally row 28)
# has shifted up to row 27. Delete it too.
$ws.Rows(27).Delete()

# Fix up the imputed values in column C for the remaining rows.
$ws.Range("C27").Value = 10      # SC 101
$ws.Range("C28").ClearContents() # SC 105
$ws.Range("C29").ClearContents() # SC 119
$ws.Range("C30").Value = 11.4    # SC 120
$ws.Range("C32").ClearContents() # SC 193
